# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell address -> new numeric value
$updates = @{
    "K2"  = 2
    "Q2"  = 2.3
    "R2"  = 1.62
    "U2"  = 1.92
    "V2"  = 1.77
    "X2"  = 11
    "AC2" = 8
    "AD2" = 6.5
    "AG2" = 501
    "AL2" = 34
    "AV2" = 67
    "AZ2" = 67
    "BA2" = 101

    "G3"  = 2.3
    "I3"  = 3.5
    "L3"  = 4
    "M3"  = 1.05
    "N3"  = 7.5
    "U3"  = 1.92
    "V3"  = 1.77
    "X3"  = 10
    "Y3"  = 9.5
    "Z3"  = 21
    "AJ3" = 13
    "AQ3" = 41
    "AR3" = 67

    "G4"  = 1.85
    "H4"  = 3.4
    "I4"  = 4.5
    "J4"  = 2.6
    "L4"  = 5
    "M4"  = 1.08
    "O4"  = 1.4
    "AG4" = 501
    "AO4" = 10
    "AZ4" = 101

    "M5" = 1.11
    "N5" = 6.5
    "V5" = 1.67

    "I8" = 1.72
    "Q8" = 1.93
    "R8" = 1.93
    "U8" = 1.8
    "V8" = 1.91

    "M23" = 1.05
    "O23" = 1.41
    "P23" = 2.62
    "R23" = 1.57
    "V23" = 1.69

    "R24" = 1.6

    "M38" = 1.08
    "O38" = 1.44
    "P38" = 2.63
    "R38" = 1.53

    "M39" = 1.05
    "O39" = 1.29
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
